$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9169.333000000001
$ws.Range("I9").Value = 14504.714
$ws.Range("J9").Value = 1699.8
$ws.Range("K9").Value = 14504.714
$ws.Range("L9").Value = 1699.8
$ws.Range("M9").Value = -14335.714
$ws.Range("N9").Value = -2037.8
$ws.Range("H18").Value = 1237.5
$ws.Range("J18").Value = 1175
$ws.Range("L18").Value = 1175
$ws.Range("N18").Value = -1743
$ws.Range("H40").Value = 4148.8
$ws.Range("I40").Value = 3840.842
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 3840.842
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -3665.842
$ws.Range("N40").Value = -10350
$ws.Range("H54").Value = 31510.857
$ws.Range("I54").Value = 25144
$ws.Range("K54").Value = 25144
$ws.Range("M54").Value = -24658
$ws.Range("H55").Value = 231.22223
$ws.Range("J55").Value = 583.3333
$ws.Range("L55").Value = 583.3333
$ws.Range("N55").Value = -1011.3333
$ws.Range("H80").Value = 898.3333
$ws.Range("I80").Value = 947.5
$ws.Range("J80").Value = 800
$ws.Range("K80").Value = 2842.5
$ws.Range("L80").Value = 2400
$ws.Range("M80").Value = -1844.5
$ws.Range("N80").Value = -4396
$ws.Range("H83").Value = 898.3333
$ws.Range("I83").Value = 947.5
$ws.Range("J83").Value = 800
$ws.Range("K83").Value = 8527.5
$ws.Range("L83").Value = 7200
$ws.Range("M83").Value = -3535.5
$ws.Range("N83").Value = -17184
$ws.Range("H98").Value = 1501.7646
$ws.Range("I98").Value = 1634
$ws.Range("J98").Value = 884.6667
$ws.Range("K98").Value = 1634
$ws.Range("L98").Value = 884.6667
$ws.Range("M98").Value = -136
$ws.Range("N98").Value = -3880.6667
$ws.Range("H122").Value = 1501.7646
$ws.Range("I122").Value = 1634
$ws.Range("J122").Value = 884.6667
$ws.Range("K122").Value = 4902
$ws.Range("L122").Value = 2654.0001
$ws.Range("M122").Value = -2452
$ws.Range("N122").Value = -7554.0001
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H129").Value = 2382.1428
$ws.Range("J129").Value = 3498.5
$ws.Range("L129").Value = 10495.5
$ws.Range("N129").Value = -20495.5
$ws.Range("H133").Value = 77812.5
$ws.Range("J133").Value = 77812.5
$ws.Range("L133").Value = 77812.5
$ws.Range("N133").Value = -87932.5
$ws.Range("H135").Value = 125000710
$ws.Range("I135").Value = 125000710
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 1125006390
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1125003855
$ws.Range("H137").Value = 3331.6785
$ws.Range("I137").Value = 2909.4375
$ws.Range("K137").Value = 8728.3125
$ws.Range("M137").Value = -6178.3125
$ws.Range("H138").Value = 2307.8071
$ws.Range("I138").Value = 2221.524
$ws.Range("J138").Value = 2358.139
$ws.Range("K138").Value = 6664.572
$ws.Range("L138").Value = 7074.417
$ws.Range("M138").Value = -1524.572
$ws.Range("N138").Value = -17354.417
$ws.Range("H141").Value = 2718.625
$ws.Range("I141").Value = 2718.625
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8155.875
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2975.875
$ws.Range("N124","N135","N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 551.6667
$ws.Range("I5").Value = 303.33334
$ws.Range("K5").Value = 303.33334
$ws.Range("M5").Value = -191.33334
$ws.Range("H7").Value = 100712
$ws.Range("J7").Value = 100712
$ws.Range("L7").Value = 100712
$ws.Range("N7").Value = -100940
$ws.Range("H32").Value = 5733.636
$ws.Range("I32").Value = 3012.2144
$ws.Range("J32").Value = 10496.125
$ws.Range("K32").Value = 3012.2144
$ws.Range("L32").Value = 10496.125
$ws.Range("M32").Value = -2725.2144
$ws.Range("N32").Value = -11070.125
$ws.Range("H74").Value = 83341740
$ws.Range("I74").Value = 90917816
$ws.Range("K74").Value = 90917816
$ws.Range("M74").Value = -90916942
$ws.Range("H77").Value = 83341740
$ws.Range("I77").Value = 90917816
$ws.Range("K77").Value = 454589080
$ws.Range("M77").Value = -454584712
$ws.Range("H110").Value = 74768.42999999999
$ws.Range("I110").Value = 113750.78
$ws.Range("K110").Value = 113750.78
$ws.Range("M110").Value = -111705.78
$ws.Range("H122").Value = 2000.6
$ws.Range("I122").Value = 1500.6666
$ws.Range("K122").Value = 4501.9998
$ws.Range("M122").Value = -2051.9998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 551.6667
$ws.Range("I4").Value = 303.33334
$ws.Range("K4").Value = 303.33334
$ws.Range("M4").Value = -188.33334
$ws.Range("H86").Value = 2064.3872
$ws.Range("I86").Value = 2173.182
$ws.Range("J86").Value = 1798.4445
$ws.Range("K86").Value = 2173.182
$ws.Range("L86").Value = 1798.4445
$ws.Range("M86").Value = -1050.182
$ws.Range("N86").Value = -4044.4445
$ws.Range("H89").Value = 2064.3872
$ws.Range("I89").Value = 2173.182
$ws.Range("J89").Value = 1798.4445
$ws.Range("K89").Value = 10865.91
$ws.Range("L89").Value = 8992.2225
$ws.Range("M89").Value = -5249.91
$ws.Range("N89").Value = -20224.2225
$ws.Range("H105").Value = 3052.1072
$ws.Range("I105").Value = 2727.238
$ws.Range("K105").Value = 2727.238
$ws.Range("M105").Value = -980.2379999999998
$ws.Range("H107").Value = 102067.9
$ws.Range("I107").Value = 1531.6666
$ws.Range("J107").Value = 252872.25
$ws.Range("K107").Value = 1531.6666
$ws.Range("L107").Value = 252872.25
$ws.Range("M107").Value = 388.3334
$ws.Range("N107").Value = -256712.25
$ws.Range("H134").Value = 13160176
$ws.Range("I134").Value = 14287564
$ws.Range("J134").Value = 7320
$ws.Range("K134").Value = 42862692
$ws.Range("L134").Value = 21960
$ws.Range("M134").Value = -42860157
$ws.Range("N134").Value = -27030

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 748.5714
$ws.Range("I19").Value = 285
$ws.Range("K19").Value = 285
$ws.Range("M19").Value = -115
$ws.Range("H23").Value = 109
$ws.Range("I23").Value = 109
$ws.Range("K23").Value = 109
$ws.Range("M23").Value = 131
$ws.Range("H24").Value = 748.5714
$ws.Range("I24").Value = 285
$ws.Range("K24").Value = 285
$ws.Range("M24").Value = -115
$ws.Range("H27").Value = 109
$ws.Range("I27").Value = 109
$ws.Range("K27").Value = 109
$ws.Range("M27").Value = 83
$ws.Range("H62").Value = 13166.667
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 17250
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 17250
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -18498
$ws.Range("H65").Value = 13166.667
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 17250
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 86250
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -92490
$ws.Range("H99").Value = 2128
$ws.Range("I99").Value = 1547
$ws.Range("K99").Value = 1547
$ws.Range("M99").Value = -49
$ws.Range("H122").Value = 3000.5217
$ws.Range("I122").Value = 2906.1
$ws.Range("K122").Value = 8718.299999999999
$ws.Range("M122").Value = -6268.299999999999
$ws.Range("H126").Value = 2128
$ws.Range("I126").Value = 1547
$ws.Range("K126").Value = 4641
$ws.Range("M126").Value = -2171
$ws.Range("H141").Value = 541250
$ws.Range("J141").Value = 688333.3
$ws.Range("L141").Value = 688333.3
$ws.Range("N141").Value = -698693.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 182.42857
$ws.Range("I2").Value = 202.5
$ws.Range("J2").Value = 174.4
$ws.Range("K2").Value = 1215
$ws.Range("L2").Value = 1046.4
$ws.Range("M2").Value = -1102
$ws.Range("N2").Value = -1272.4
$ws.Range("H14").Value = 377.25
$ws.Range("I14").Value = 377.25
$ws.Range("K14").Value = 1131.75
$ws.Range("M14").Value = -958.75
$ws.Range("H15").Value = 113.333336
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("H37").Value = 138100
$ws.Range("J37").Value = 138100
$ws.Range("L37").Value = 414300
$ws.Range("N37").Value = -414524
$ws.Range("H80").Value = 3494.5
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 9000
$ws.Range("N80").Value = -10872
$ws.Range("H83").Value = 3494.5
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 27000
$ws.Range("N83").Value = -36360
$ws.Range("H107").Value = 1414.2106
$ws.Range("I107").Value = 558
$ws.Range("J107").Value = 1809.3846
$ws.Range("K107").Value = 1674
$ws.Range("L107").Value = 5428.1538
$ws.Range("M107").Value = 246
$ws.Range("N107").Value = -9268.1538
$ws.Range("H131").Value = 1995.3125
$ws.Range("I131").Value = 2132.7144
$ws.Range("K131").Value = 6398.1432
$ws.Range("M131").Value = -1358.1432
$ws.Range("H137").Value = 7144159
$ws.Range("I137").Value = 7144159
$ws.Range("K137").Value = 21432477
$ws.Range("M137").Value = -21427377
$ws.Range("N15").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 40000
$ws.Range("I35").Value = 37500
$ws.Range("J35").Value = 45000
$ws.Range("K35").Value = 37500
$ws.Range("L35").Value = 45000
$ws.Range("M35").Value = -37202
$ws.Range("N35").Value = -45596
$ws.Range("H122").Value = 4572.4287
$ws.Range("I122").Value = 501.75
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 1505.25
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = 944.75
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 5004706
$ws.Range("I132").Value = 5955790.5
$ws.Range("K132").Value = 17867371.5
$ws.Range("M132").Value = -17864841.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1759
$ws.Range("I16").Value = 450.05264
$ws.Range("J16").Value = 3535.4285
$ws.Range("K16").Value = 450.05264
$ws.Range("L16").Value = 3535.4285
$ws.Range("M16").Value = -280.05264
$ws.Range("N16").Value = -3875.4285
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H46").Value = 874
$ws.Range("I46").Value = 874
$ws.Range("K46").Value = 874
$ws.Range("M46").Value = -686
$ws.Range("H55").Value = 183.38889
$ws.Range("I55").Value = 131.89473
$ws.Range("K55").Value = 131.89473
$ws.Range("M55").Value = 41.10526999999999
$ws.Range("H122").Value = 5689.7095
$ws.Range("I122").Value = 5568.3105
$ws.Range("K122").Value = 16704.9315
$ws.Range("M122").Value = -14254.9315
$ws.Range("M18","N18").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 5000
$ws.Range("J11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("N11").Value = -5284
$ws.Range("H37").Value = 29
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 29
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 29
$ws.Range("N37").Value = -435
$ws.Range("H81").Value = 3010.6155
$ws.Range("I81").Value = 3111.5
$ws.Range("K81").Value = 6223
$ws.Range("M81").Value = -5162
$ws.Range("H84").Value = 3010.6155
$ws.Range("I84").Value = 3111.5
$ws.Range("K84").Value = 31115
$ws.Range("M84").Value = -25811
$ws.Range("H86").Value = 40325
$ws.Range("J86").Value = 40325
$ws.Range("L86").Value = 40325
$ws.Range("N86").Value = -42571
$ws.Range("H89").Value = 40325
$ws.Range("J89").Value = 40325
$ws.Range("L89").Value = 201625
$ws.Range("N89").Value = -212857
$ws.Range("H107").Value = 716.8077
$ws.Range("I107").Value = 462.14285
$ws.Range("K107").Value = 1386.42855
$ws.Range("M107").Value = 533.5714499999999
$ws.Range("H126").Value = 1979.6364
$ws.Range("I126").Value = 1937.125
$ws.Range("K126").Value = 5811.375
$ws.Range("M126").Value = -3341.375
$ws.Range("M37").ClearContents()
